$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.852.01"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.088.31"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +5.08%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.94"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.76"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.48%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.082.56"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("E12").Value = "  +5.72%  "
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.48"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.599.45"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.881.62"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("E18").Value = "  +4.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.086.78"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +10.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "465.66"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.76%  "
$ws.Range("E22").Value = "  +4.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.52"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.38"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("E25").Value = "  +7.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.98"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.16"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.98"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.26"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.45%  "
$ws.Range("E34").Value = "  +4.09%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.66%  "
$ws.Range("E37").Value = "  +3.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.13"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.11"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.319"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.22"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("E42").Value = "  +1.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.68"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.91%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "384.76"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.772.25"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.09"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.82"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +7.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.36%  "
